# ProviderSubmissionsReportTemplate.xlsx edit
# Commit: "Updated templates to have different worksheet name to correct
# tests. Renamed the model lists in data quality to improve readability
# of the code."
#
# The only user-facing change inside this template workbook is the
# worksheet's display name: "Data Quality" -> "Provider Submissions"
# (the "model list" renames mentioned in the commit message are C#
# code-level changes in the report-service project, not part of this
# spreadsheet). We also carry over the last active-cell selection that
# was recorded the last time the template was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab/displayed name.
$ws.Name = "Provider Submissions"

# Match the recorded selection in the saved template.
$ws.Activate()
$ws.Range("B32").Select()
